# Added links to course list
# Populate column E ("course code") with the DTU course numbers that
# correspond to each course name already listed in column A, formatted
# as Text (numFmtId 49 / "@") so that leading zeros are preserved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# (Row, Code) pairs, given in the exact order the values were originally
# entered (row 2, then row 3, then row 1, then sequentially from row 4
# onward; row 20 - the exchange-semester row - has no course code).
$entries = @(
    [PSCustomObject]@{ Row=2;  Code="01017" },
    [PSCustomObject]@{ Row=3;  Code="02631" },
    [PSCustomObject]@{ Row=1;  Code="02461" },
    [PSCustomObject]@{ Row=4;  Code="02403" },
    [PSCustomObject]@{ Row=5;  Code="01005" },
    [PSCustomObject]@{ Row=6;  Code="02462" },
    [PSCustomObject]@{ Row=7;  Code="02464" },
    [PSCustomObject]@{ Row=8;  Code="02105" },
    [PSCustomObject]@{ Row=9;  Code="26028" },
    [PSCustomObject]@{ Row=10; Code="02809" },
    [PSCustomObject]@{ Row=11; Code="02445" },
    [PSCustomObject]@{ Row=12; Code="01035" },
    [PSCustomObject]@{ Row=13; Code="02450" },
    [PSCustomObject]@{ Row=14; Code="10024" },
    [PSCustomObject]@{ Row=15; Code="34315" },
    [PSCustomObject]@{ Row=16; Code="02463" },
    [PSCustomObject]@{ Row=17; Code="02466" },
    [PSCustomObject]@{ Row=18; Code="02465" },
    [PSCustomObject]@{ Row=19; Code="42611" },
    [PSCustomObject]@{ Row=21; Code="02182" },
    [PSCustomObject]@{ Row=22; Code="02502" },
    [PSCustomObject]@{ Row=23; Code="02170" }
)

foreach ($entry in $entries) {
    $cell = $ws.Range("E$($entry.Row)")
    # Format as Text first so the leading zeros in the course numbers
    # ("01017", "02461", ...) are not stripped by automatic number
    # conversion.
    $cell.NumberFormat = "@"
    $cell.Value = $entry.Code
}

# Match the new selection left behind in the saved workbook.
$ws.Range("C32").Select() | Out-Null
